$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 224.74074
$ws.Range("I53").Value = 88.47059
$ws.Range("K53").Value = 88.47059
$ws.Range("M53").Value = 548.52941
# Row 106
$ws.Range("H106").Value = 2453.75
$ws.Range("I106").Value = 2672.5
$ws.Range("J106").Value = 1797.5
$ws.Range("K106").Value = 2672.5
$ws.Range("L106").Value = 1797.5
$ws.Range("M106").Value = -2041.5
$ws.Range("N106").Value = -3059.5
# Row 113
$ws.Range("H113").Value = 2727.7917
$ws.Range("I113").Value = 3044.3333
$ws.Range("J113").Value = 2200.2222
$ws.Range("K113").Value = 3044.3333
$ws.Range("L113").Value = 2200.2222
$ws.Range("M113").Value = 209.6667000000002
$ws.Range("N113").Value = -8708.2222

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2198076
$ws.Range("I32").Value = 4165.8774
$ws.Range("J32").Value = 15635776
$ws.Range("K32").Value = 4165.8774
$ws.Range("L32").Value = 15635776
$ws.Range("M32").Value = -3878.8774
$ws.Range("N32").Value = -15636350
# Row 61
$ws.Range("H61").Value = 1518.6666
$ws.Range("I61").Value = 1379.7391
$ws.Range("K61").Value = 1379.7391
$ws.Range("M61").Value = -1167.7391
# Row 110
$ws.Range("H110").Value = 631.8889
$ws.Range("I110").Value = 648.375
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 648.375
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 1396.625
$ws.Range("N110").Value = -4590
# Row 133
$ws.Range("H133").Value = 33800
$ws.Range("J133").Value = 33800
$ws.Range("L133").Value = 33800
$ws.Range("N133").Value = -38860
# Row 136
$ws.Range("H136").Value = 1518.6666
$ws.Range("I136").Value = 1379.7391
$ws.Range("K136").Value = 4139.2173
$ws.Range("M136").Value = -1589.2173
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1324.2222
$ws.Range("I94").Value = 859.7143
$ws.Range("K94").Value = 859.7143
$ws.Range("M94").Value = -408.7143
# Row 107
$ws.Range("H107").Value = 11885.167
$ws.Range("I107").Value = 867.75
$ws.Range("K107").Value = 867.75
$ws.Range("M107").Value = 1052.25
# Row 132
$ws.Range("H132").Value = 194285.72
$ws.Range("J132").Value = 194285.72
$ws.Range("L132").Value = 194285.72
$ws.Range("N132").Value = -204405.72

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 952.8333
$ws.Range("I16").Value = 923.4
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 923.4
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -636.4
$ws.Range("N16").Value = -1674
# Row 58
$ws.Range("H58").Value = 6228.1
$ws.Range("I58").Value = 930.2857
$ws.Range("J58").Value = 9080.77
$ws.Range("K58").Value = 930.2857
$ws.Range("L58").Value = 9080.77
$ws.Range("M58").Value = -727.2857
$ws.Range("N58").Value = -9486.77
# Row 113
$ws.Range("H113").Value = 952.8333
$ws.Range("I113").Value = 923.4
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 923.4
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1246.6
$ws.Range("N113").Value = -5440
# Row 132
$ws.Range("H132").Value = 3154.9524
$ws.Range("I132").Value = 2379.9167
$ws.Range("J132").Value = 4188.3335
$ws.Range("K132").Value = 7139.750100000001
$ws.Range("L132").Value = 12565.0005
$ws.Range("M132").Value = -4609.750100000001
$ws.Range("N132").Value = -17625.0005
# Row 136
$ws.Range("H136").Value = 6228.1
$ws.Range("I136").Value = 930.2857
$ws.Range("J136").Value = 9080.77
$ws.Range("K136").Value = 2790.8571
$ws.Range("L136").Value = 27242.31
$ws.Range("M136").Value = -240.8571000000002
$ws.Range("N136").Value = -32342.31
# Row 140
$ws.Range("H140").Value = 54244
$ws.Range("J140").Value = 54244
$ws.Range("L140").Value = 54244
$ws.Range("N140").Value = -64604

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1472.5862
$ws.Range("I5").Value = 1131.2
$ws.Range("J5").Value = 3606.25
$ws.Range("K5").Value = 3393.6
$ws.Range("L5").Value = 10818.75
$ws.Range("M5").Value = -3281.6
$ws.Range("N5").Value = -11042.75
# Row 107
$ws.Range("H107").Value = 707840.5600000001
$ws.Range("I107").Value = 846.9
$ws.Range("J107").Value = 7777777
$ws.Range("K107").Value = 2540.7
$ws.Range("L107").Value = 23333331
$ws.Range("M107").Value = -620.6999999999998
$ws.Range("N107").Value = -23337171
# Row 117
$ws.Range("H117").Value = 2501.1538
$ws.Range("I117").Value = 280
$ws.Range("J117").Value = 2686.25
$ws.Range("K117").Value = 840
$ws.Range("L117").Value = 8058.75
$ws.Range("M117").Value = 2602
$ws.Range("N117").Value = -14942.75
# Row 129
$ws.Range("H129").Value = 2545.1738
$ws.Range("I129").Value = 1634.875
$ws.Range("J129").Value = 3030.6667
$ws.Range("K129").Value = 4904.625
$ws.Range("L129").Value = 9092.000100000001
$ws.Range("M129").Value = 95.375
$ws.Range("N129").Value = -19092.0001
# Row 135
$ws.Range("H135").Value = 1472.5862
$ws.Range("I135").Value = 1131.2
$ws.Range("J135").Value = 3606.25
$ws.Range("K135").Value = 10180.8
$ws.Range("L135").Value = 32456.25
$ws.Range("M135").Value = -7645.800000000001
$ws.Range("N135").Value = -37526.25

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4612.25
$ws.Range("I70").Value = 4099.923
$ws.Range("J70").Value = 5563.7144
$ws.Range("K70").Value = 4099.923
$ws.Range("L70").Value = 5563.7144
$ws.Range("M70").Value = -3829.923
$ws.Range("N70").Value = -6103.7144
# Row 73
$ws.Range("H73").Value = 4612.25
$ws.Range("I73").Value = 4099.923
$ws.Range("J73").Value = 5563.7144
$ws.Range("K73").Value = 4099.923
$ws.Range("L73").Value = 5563.7144
$ws.Range("M73").Value = -3163.923
$ws.Range("N73").Value = -7435.7144
# Row 113
$ws.Range("H113").Value = 1175
$ws.Range("I113").Value = 1235.3334
$ws.Range("J113").Value = 813
$ws.Range("K113").Value = 1235.3334
$ws.Range("L113").Value = 813
$ws.Range("M113").Value = 934.6666
$ws.Range("N113").Value = -5153
# Row 122
$ws.Range("H122").Value = 2406.0605
$ws.Range("J122").Value = 3346.3845
$ws.Range("L122").Value = 10039.1535
$ws.Range("N122").Value = -14939.1535

$ws = $wb.Worksheets.Item("LTW")
# Row 135
$ws.Range("H135").Value = 57036.555
$ws.Range("J135").Value = 57036.555
$ws.Range("L135").Value = 57036.555
$ws.Range("N135").Value = -67176.55499999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 813.5333000000001
$ws.Range("I107").Value = 865.86957
$ws.Range("J107").Value = 641.5714
$ws.Range("K107").Value = 2597.60871
$ws.Range("L107").Value = 1924.7142
$ws.Range("M107").Value = -677.60871
$ws.Range("N107").Value = -5764.7142
# Row 132
$ws.Range("H132").Value = 1858.4736
$ws.Range("I132").Value = 1534.1333
$ws.Range("J132").Value = 3074.75
$ws.Range("K132").Value = 4602.3999
$ws.Range("L132").Value = 9224.25
$ws.Range("M132").Value = -2072.3999
$ws.Range("N132").Value = -14284.25
# Row 136
$ws.Range("H136").Value = 1592.4193
$ws.Range("I136").Value = 1528.8334
$ws.Range("K136").Value = 4586.5002
$ws.Range("M136").Value = -2036.5002
